$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date-like text to be stored as plain text (matching the source
# data, which uses inline/shared strings like "MM/DD/YYYY" rather than a
# real date value), then drop back to the default "Normal" style so no
# number-format is left applied to the cell.
$ws.Range("A27").NumberFormat = "@"
$ws.Range("A27").Value = "09/28/2025"
$ws.Range("A27").Style = "Normal"

$ws.Range("B27").Value = 0.1383344084454727
$ws.Range("C27").Value = 0.8616655915545273
